# TCER Sheet.xlsx - "Transaction Service" -> add "Transaction Dao" section
# Commit message: "transaction service dao tcer done"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Transaction Service" tab

# -----------------------------------------------------------------
# 1) Row 4 (R001/T001 saveTransaction comment) - reword Comment text
# -----------------------------------------------------------------
$ws.Range("F4").Value = "if it  transaction is added "

# -----------------------------------------------------------------
# 2) New row 5 - negative case for saveTransaction, inserted right
#    below row 4 (same Requirement/TestCase block)
# -----------------------------------------------------------------
$ws.Range("A5").Value = "R001"
$ws.Range("B5").Value = "T002"
$ws.Range("C5").Value = "saveTransaction(Transaction transaction)"
$ws.Range("D5").Value = $false
$ws.Range("F5").Value = "if it  transaction is not added "

# -----------------------------------------------------------------
# 3) Row 10 (getMaxTransactionId) - updated Expected Result + new Comment
# -----------------------------------------------------------------
$ws.Range("D10").Value = "latest transaction id"
$ws.Range("F10").Value = "to get recent transaction id"

# -----------------------------------------------------------------
# 4) New section header (row 19, merged A19:F19) - "Transaction Dao"
#    styled like the existing section header in row 2.
# -----------------------------------------------------------------
$ws.Range("A19:F19").Merge()
$ws.Range("A19").Value = "Transaction Dao"
$sectionFont = $ws.Range("A2").Font
$newSectionRange = $ws.Range("A19:F19")
$newSectionRange.Font.Bold = $sectionFont.Bold
$newSectionRange.Font.Size = $sectionFont.Size
$newSectionRange.HorizontalAlignment = $ws.Range("A2").HorizontalAlignment
$ws.Rows.Item(19).RowHeight = 18

# -----------------------------------------------------------------
# 5) New column-header row (row 20) - same headers as row 3.
# -----------------------------------------------------------------
$ws.Range("A20").Value = "Requirement #"
$ws.Range("B20").Value = "TestCaseId"
$ws.Range("C20").Value = "TestCase"
$ws.Range("D20").Value = "Expected Result"
$ws.Range("E20").Value = "Actual Result"
$ws.Range("F20").Value = "Comment"
$headerRange = $ws.Range("A20:F20")
$headerFont = $ws.Range("A3").Font
$headerRange.Font.Bold = $headerFont.Bold
$headerRange.Font.Size = $headerFont.Size
$headerRange.Borders.LineStyle = $ws.Range("A3").Borders.LineStyle
$ws.Rows.Item(20).RowHeight = 15.6

# -----------------------------------------------------------------
# 6) Row 21/22 - saveTransactionByUserName(String userName): R001
# -----------------------------------------------------------------
$ws.Range("A21").Value = "R001"
$ws.Range("B21").Value = "T001"
$ws.Range("C21").Value = "saveTransactionByUserName(String userName)"
$ws.Range("D21").Value = $true
$ws.Range("F21").Value = "If transaction is saved"

$ws.Range("A22").Value = "R001"
$ws.Range("B22").Value = "T002"
$ws.Range("C22").Value = "saveTransactionByUserName(String userName)"
$ws.Range("D22").Value = $false
$ws.Range("F22").Value = "if transaction is not saved"

# -----------------------------------------------------------------
# 7) Row 24 - getMaxTransaction(): R002
# -----------------------------------------------------------------
$ws.Range("A24").Value = "R002"
$ws.Range("B24").Value = "T001"
$ws.Range("C24").Value = "getMaxTransaction()"
$ws.Range("D24").Value = "latest Transaction id"
$ws.Range("F24").Value = "to get recent transaction id"

# -----------------------------------------------------------------
# 8) Row 26/27 - findAll(): R003
# -----------------------------------------------------------------
$ws.Range("A26").Value = "R003"
$ws.Range("B26").Value = "T001"
$ws.Range("C26").Value = "findAll()"
$ws.Range("D26").Value = "List All transactions"
$ws.Range("F26").Value = "if there is atleast one transaction done by using application"

$ws.Range("A27").Value = "R003"
$ws.Range("B27").Value = "T002"
$ws.Range("C27").Value = "findAll()"
$ws.Range("D27").Value = "null"
$ws.Range("F27").Value = "if there is transaction done by using application"

# -----------------------------------------------------------------
# 9) Row 29/30 - saveTransactionByUserName(String userName) again: R004
#    (list transactions of a particular user)
# -----------------------------------------------------------------
$ws.Range("A29").Value = "R004"
$ws.Range("B29").Value = "T001"
$ws.Range("C29").Value = "saveTransactionByUserName(String userName)"
$ws.Range("D29").Value = "List transactions of particular user"
$ws.Range("F29").Value = "if atleast one transaction exists for customer"

$ws.Range("A30").Value = "R004"
$ws.Range("B30").Value = "T002"
$ws.Range("C30").Value = "saveTransactionByUserName(String userName)"
$ws.Range("D30").Value = "null"
$ws.Range("F30").Value = "if no transaction is done for customer"

# -----------------------------------------------------------------
# 10) Column widths / view tweaks to roughly match the refreshed layout
# -----------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 40.1
$ws.Columns.Item(4).ColumnWidth = 40.1
$ws.Columns.Item(6).ColumnWidth = 49.7

$ws.Range("D31").Select()
